# Apply refreshed TPM-derived NATMI metrics to Sheet1.
# Only the numeric result columns (G-T, excluding the constant
# detection-rate/count columns) change between the old and new TPM run;
# ligand/receptor identity columns (A-F, K, L) are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 19.163986
$ws.Range("H2").Value = 57.491958
$ws.Range("I2").Value = 0.1197574615923936
$ws.Range("J2").Value = 0.1197574615923936
$ws.Range("M2").Value = 51.15371566666666
$ws.Range("N2").Value = 153.461147
$ws.Range("O2").Value = 0.3311207986511828
$ws.Range("P2").Value = 0.3311207986511828
$ws.Range("Q2").Value = 980.3090908839806
$ws.Range("R2").Value = 8822.781817955825
$ws.Range("S2").Value = 0.03965418632691172
$ws.Range("T2").Value = 0.03965418632691172
$ws.Range("G3").Value = 19.163986
$ws.Range("H3").Value = 57.491958
$ws.Range("I3").Value = 0.1197574615923936
$ws.Range("J3").Value = 0.1197574615923936
$ws.Range("M3").Value = 53.36146666666667
$ws.Range("O3").Value = 0.3454116915964105
$ws.Range("P3").Value = 0.3454116915964106
$ws.Range("Q3").Value = 1022.618400139467
$ws.Range("R3").Value = 9203.565601255201
$ws.Range("S3").Value = 0.04136562738992083
$ws.Range("T3").Value = 0.04136562738992084
$ws.Range("G4").Value = 19.163986
$ws.Range("H4").Value = 57.491958
$ws.Range("I4").Value = 0.1197574615923936
$ws.Range("J4").Value = 0.1197574615923936
$ws.Range("M4").Value = 25.00653133333333
$ws.Range("N4").Value = 75.019594
$ws.Range("O4").Value = 0.1618686447050176
$ws.Range("P4").Value = 0.1618686447050176
$ws.Range("Q4").Value = 479.2248163805613
$ws.Range("R4").Value = 4313.023347425052
$ws.Range("S4").Value = 0.01938497800127395
$ws.Range("T4").Value = 0.01938497800127395
$ws.Range("G5").Value = 19.163986
$ws.Range("H5").Value = 57.491958
$ws.Range("I5").Value = 0.1197574615923936
$ws.Range("J5").Value = 0.1197574615923936
$ws.Range("M5").Value = 6.481347
$ws.Range("N5").Value = 19.444041
$ws.Range("O5").Value = 0.0419541135381084
$ws.Range("P5").Value = 0.0419541135381084
$ws.Range("Q5").Value = 124.208443169142
$ws.Range("R5").Value = 1117.875988522278
$ws.Range("S5").Value = 0.005024318140682936
$ws.Range("T5").Value = 0.005024318140682937
$ws.Range("G6").Value = 19.163986
$ws.Range("H6").Value = 57.491958
$ws.Range("I6").Value = 0.1197574615923936
$ws.Range("J6").Value = 0.1197574615923936
$ws.Range("M6").Value = 18.483507
$ws.Range("N6").Value = 55.450521
$ws.Range("O6").Value = 0.1196447515092806
$ws.Range("P6").Value = 0.1196447515092806
$ws.Range("Q6").Value = 354.217669378902
$ws.Range("R6").Value = 3187.959024410118
$ws.Range("S6").Value = 0.01432835173360415
$ws.Range("T6").Value = 0.01432835173360415
$ws.Range("I7").Value = 0.150345281456851
$ws.Range("J7").Value = 0.1503452814568511
$ws.Range("M7").Value = 51.15371566666666
$ws.Range("N7").Value = 153.461147
$ws.Range("O7").Value = 0.3311207986511828
$ws.Range("P7").Value = 0.3311207986511828
$ws.Range("Q7").Value = 1230.694473846655
$ws.Range("R7").Value = 11076.2502646199
$ws.Range("S7").Value = 0.04978244966942938
$ws.Range("T7").Value = 0.04978244966942939
$ws.Range("I8").Value = 0.150345281456851
$ws.Range("J8").Value = 0.1503452814568511
$ws.Range("M8").Value = 53.36146666666667
$ws.Range("O8").Value = 0.3454116915964105
$ws.Range("P8").Value = 0.3454116915964106
$ws.Range("S8").Value = 0.05193101799154936
$ws.Range("T8").Value = 0.05193101799154938
$ws.Range("I9").Value = 0.150345281456851
$ws.Range("J9").Value = 0.1503452814568511
$ws.Range("M9").Value = 25.00653133333333
$ws.Range("N9").Value = 75.019594
$ws.Range("O9").Value = 0.1618686447050176
$ws.Range("P9").Value = 0.1618686447050176
$ws.Range("Q9").Value = 601.6258940513437
$ws.Range("R9").Value = 5414.633046462093
$ws.Range("S9").Value = 0.0243361869472149
$ws.Range("T9").Value = 0.0243361869472149
$ws.Range("I10").Value = 0.150345281456851
$ws.Range("J10").Value = 0.1503452814568511
$ws.Range("M10").Value = 6.481347
$ws.Range("N10").Value = 19.444041
$ws.Range("O10").Value = 0.0419541135381084
$ws.Range("P10").Value = 0.0419541135381084
$ws.Range("Q10").Value = 155.933109296699
$ws.Range("R10").Value = 1403.397983670291
$ws.Range("S10").Value = 0.006307603008159591
$ws.Range("T10").Value = 0.006307603008159592
$ws.Range("I11").Value = 0.150345281456851
$ws.Range("J11").Value = 0.1503452814568511
$ws.Range("M11").Value = 18.483507
$ws.Range("N11").Value = 55.450521
$ws.Range("O11").Value = 0.1196447515092806
$ws.Range("P11").Value = 0.1196447515092806
$ws.Range("Q11").Value = 444.690080197419
$ws.Range("R11").Value = 4002.210721776771
$ws.Range("S11").Value = 0.0179880238404978
$ws.Range("T11").Value = 0.0179880238404978
$ws.Range("G12").Value = 61.341815
$ws.Range("H12").Value = 184.025445
$ws.Range("I12").Value = 0.3833304853108436
$ws.Range("J12").Value = 0.3833304853108436
$ws.Range("M12").Value = 51.15371566666666
$ws.Range("N12").Value = 153.461147
$ws.Range("O12").Value = 0.3311207986511828
$ws.Range("P12").Value = 0.3311207986511828
$ws.Range("Q12").Value = 3137.861762987268
$ws.Range("R12").Value = 28240.75586688541
$ws.Range("S12").Value = 0.126928696443472
$ws.Range("T12").Value = 0.126928696443472
$ws.Range("G13").Value = 61.341815
$ws.Range("H13").Value = 184.025445
$ws.Range("I13").Value = 0.3833304853108436
$ws.Range("J13").Value = 0.3833304853108436
$ws.Range("M13").Value = 53.36146666666667
$ws.Range("O13").Value = 0.3454116915964105
$ws.Range("P13").Value = 0.3454116915964106
$ws.Range("Q13").Value = 3273.289216395333
$ws.Range("R13").Value = 29459.602947558
$ws.Range("S13").Value = 0.1324068313716915
$ws.Range("T13").Value = 0.1324068313716915
$ws.Range("G14").Value = 61.341815
$ws.Range("H14").Value = 184.025445
$ws.Range("I14").Value = 0.3833304853108436
$ws.Range("J14").Value = 0.3833304853108436
$ws.Range("M14").Value = 25.00653133333333
$ws.Range("N14").Value = 75.019594
$ws.Range("O14").Value = 0.1618686447050176
$ws.Range("P14").Value = 0.1618686447050176
$ws.Range("Q14").Value = 1533.946018841036
$ws.Range("R14").Value = 13805.51416956933
$ws.Range("S14").Value = 0.06204918613138292
$ws.Range("T14").Value = 0.06204918613138292
$ws.Range("G15").Value = 61.341815
$ws.Range("H15").Value = 184.025445
$ws.Range("I15").Value = 0.3833304853108436
$ws.Range("J15").Value = 0.3833304853108436
$ws.Range("M15").Value = 6.481347
$ws.Range("N15").Value = 19.444041
$ws.Range("O15").Value = 0.0419541135381084
$ws.Range("P15").Value = 0.0419541135381084
$ws.Range("Q15").Value = 397.5775886248049
$ws.Range("R15").Value = 3578.198297623244
$ws.Range("S15").Value = 0.01608229070334933
$ws.Range("T15").Value = 0.01608229070334933
$ws.Range("G16").Value = 61.341815
$ws.Range("H16").Value = 184.025445
$ws.Range("I16").Value = 0.3833304853108436
$ws.Range("J16").Value = 0.3833304853108436
$ws.Range("M16").Value = 18.483507
$ws.Range("N16").Value = 55.450521
$ws.Range("O16").Value = 0.1196447515092806
$ws.Range("P16").Value = 0.1196447515092806
$ws.Range("Q16").Value = 1133.811866945205
$ws.Range("R16").Value = 10204.30680250684
$ws.Range("S16").Value = 0.04586348066094782
$ws.Range("T16").Value = 0.04586348066094783
$ws.Range("G17").Value = 7.095824666666666
$ws.Range("H17").Value = 21.287474
$ws.Range("I17").Value = 0.04434244264135302
$ws.Range("J17").Value = 0.04434244264135302
$ws.Range("M17").Value = 51.15371566666666
$ws.Range("N17").Value = 153.461147
$ws.Range("O17").Value = 0.3311207986511828
$ws.Range("P17").Value = 0.3311207986511828
$ws.Range("Q17").Value = 362.9777974191864
$ws.Range("R17").Value = 3266.800176772678
$ws.Range("S17").Value = 0.01468270502154908
$ws.Range("T17").Value = 0.01468270502154908
$ws.Range("G18").Value = 7.095824666666666
$ws.Range("H18").Value = 21.287474
$ws.Range("I18").Value = 0.04434244264135302
$ws.Range("J18").Value = 0.04434244264135302
$ws.Range("M18").Value = 53.36146666666667
$ws.Range("O18").Value = 0.3454116915964105
$ws.Range("P18").Value = 0.3454116915964106
$ws.Range("Q18").Value = 378.6436114228445
$ws.Range("R18").Value = 3407.7925028056
$ws.Range("S18").Value = 0.01531639812226655
$ws.Range("T18").Value = 0.01531639812226655
$ws.Range("G19").Value = 7.095824666666666
$ws.Range("H19").Value = 21.287474
$ws.Range("I19").Value = 0.04434244264135302
$ws.Range("J19").Value = 0.04434244264135302
$ws.Range("M19").Value = 25.00653133333333
$ws.Range("N19").Value = 75.019594
$ws.Range("O19").Value = 0.1618686447050176
$ws.Range("P19").Value = 0.1618686447050176
$ws.Range("Q19").Value = 177.4419618628395
$ws.Range("R19").Value = 1596.977656765556
$ws.Range("S19").Value = 0.007177651093265796
$ws.Range("T19").Value = 0.007177651093265796
$ws.Range("G20").Value = 7.095824666666666
$ws.Range("H20").Value = 21.287474
$ws.Range("I20").Value = 0.04434244264135302
$ws.Range("J20").Value = 0.04434244264135302
$ws.Range("M20").Value = 6.481347
$ws.Range("N20").Value = 19.444041
$ws.Range("O20").Value = 0.0419541135381084
$ws.Range("P20").Value = 0.0419541135381084
$ws.Range("Q20").Value = 45.99050191582599
$ws.Range("R20").Value = 413.914517242434
$ws.Range("S20").Value = 0.001860347873132384
$ws.Range("T20").Value = 0.001860347873132384
$ws.Range("G21").Value = 7.095824666666666
$ws.Range("H21").Value = 21.287474
$ws.Range("I21").Value = 0.04434244264135302
$ws.Range("J21").Value = 0.04434244264135302
$ws.Range("M21").Value = 18.483507
$ws.Range("N21").Value = 55.450521
$ws.Range("O21").Value = 0.1196447515092806
$ws.Range("P21").Value = 0.1196447515092806
$ws.Range("Q21").Value = 131.155724897106
$ws.Range("R21").Value = 1180.401524073954
$ws.Range("S21").Value = 0.005305340531139211
$ws.Range("T21").Value = 0.005305340531139211
$ws.Range("G22").Value = 48.362939
$ws.Range("H22").Value = 145.088817
$ws.Range("I22").Value = 0.3022243289985588
$ws.Range("J22").Value = 0.3022243289985588
$ws.Range("M22").Value = 51.15371566666666
$ws.Range("N22").Value = 153.461147
$ws.Range("O22").Value = 0.3311207986511828
$ws.Range("P22").Value = 0.3311207986511828
$ws.Range("Q22").Value = 2473.944030410345
$ws.Range("R22").Value = 22265.4962736931
$ws.Range("S22").Value = 0.1000727611898206
$ws.Range("T22").Value = 0.1000727611898206
$ws.Range("G23").Value = 48.362939
$ws.Range("H23").Value = 145.088817
$ws.Range("I23").Value = 0.3022243289985588
$ws.Range("J23").Value = 0.3022243289985588
$ws.Range("M23").Value = 53.36146666666667
$ws.Range("O23").Value = 0.3454116915964105
$ws.Range("P23").Value = 0.3454116915964106
$ws.Range("Q23").Value = 2580.717357350534
$ws.Range("R23").Value = 23226.4562161548
$ws.Range("S23").Value = 0.1043918167209823
$ws.Range("T23").Value = 0.1043918167209823
$ws.Range("G24").Value = 48.362939
$ws.Range("H24").Value = 145.088817
$ws.Range("I24").Value = 0.3022243289985588
$ws.Range("J24").Value = 0.3022243289985588
$ws.Range("M24").Value = 25.00653133333333
$ws.Range("N24").Value = 75.019594
$ws.Range("O24").Value = 0.1618686447050176
$ws.Range("P24").Value = 0.1618686447050176
$ws.Range("Q24").Value = 1209.389349475589
$ws.Range("R24").Value = 10884.5041452803
$ws.Range("S24").Value = 0.04892064253188007
$ws.Range("T24").Value = 0.04892064253188007
$ws.Range("G25").Value = 48.362939
$ws.Range("H25").Value = 145.088817
$ws.Range("I25").Value = 0.3022243289985588
$ws.Range("J25").Value = 0.3022243289985588
$ws.Range("M25").Value = 6.481347
$ws.Range("N25").Value = 19.444041
$ws.Range("O25").Value = 0.0419541135381084
$ws.Range("P25").Value = 0.0419541135381084
$ws.Range("Q25").Value = 313.456989598833
$ws.Range("R25").Value = 2821.112906389497
$ws.Range("S25").Value = 0.01267955381278416
$ws.Range("T25").Value = 0.01267955381278416
$ws.Range("G26").Value = 48.362939
$ws.Range("H26").Value = 145.088817
$ws.Range("I26").Value = 0.3022243289985588
$ws.Range("J26").Value = 0.3022243289985588
$ws.Range("M26").Value = 18.483507
$ws.Range("N26").Value = 55.450521
$ws.Range("O26").Value = 0.1196447515092806
$ws.Range("P26").Value = 0.1196447515092806
$ws.Range("Q26").Value = 893.9167215470731
$ws.Range("R26").Value = 8045.250493923658
$ws.Range("S26").Value = 0.03615955474309163
$ws.Range("T26").Value = 0.03615955474309164
